# Refresh cached market-price / profit figures on several sheets.
# (plain data cells -- no formulas -- updated by the scheduled Sheets runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3288.4443
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3288.4443
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 9865.332900000001
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -10405.3329
$ws.Range("H73").Value = 3288.4443
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3288.4443
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 9865.332900000001
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -11737.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12199973
$ws.Range("I32").Value = 12504872
$ws.Range("K32").Value = 12504872
$ws.Range("M32").Value = -12504585
$ws.Range("H61").Value = 2760.5
$ws.Range("I61").Value = 2325.3264
$ws.Range("J61").Value = 5129.778
$ws.Range("K61").Value = 2325.3264
$ws.Range("L61").Value = 5129.778
$ws.Range("M61").Value = -2113.3264
$ws.Range("N61").Value = -5553.778
$ws.Range("H132").Value = 3160.7097
$ws.Range("I132").Value = 1916.5
$ws.Range("J132").Value = 6202.1113
$ws.Range("K132").Value = 5749.5
$ws.Range("L132").Value = 18606.3339
$ws.Range("M132").Value = -3219.5
$ws.Range("N132").Value = -23666.3339
$ws.Range("H136").Value = 2760.5
$ws.Range("I136").Value = 2325.3264
$ws.Range("J136").Value = 5129.778
$ws.Range("K136").Value = 6975.9792
$ws.Range("L136").Value = 15389.334
$ws.Range("M136").Value = -4425.9792
$ws.Range("N136").Value = -20489.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 725
$ws.Range("J14").Value = 725
$ws.Range("L14").Value = 725
$ws.Range("N14").Value = -1065
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H19").Value = 2143288.2
$ws.Range("I19").Value = 3000563.5
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 3000563.5
$ws.Range("L19").Value = 100
$ws.Range("M19").Value = -3000393.5
$ws.Range("N19").Value = -440
$ws.Range("H24").Value = 2143288.2
$ws.Range("I24").Value = 3000563.5
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 3000563.5
$ws.Range("L24").Value = 100
$ws.Range("M24").Value = -3000393.5
$ws.Range("N24").Value = -440
$ws.Range("H26").Value = 9000
$ws.Range("J26").Value = 9000
$ws.Range("L26").Value = 9000
$ws.Range("N26").Value = -9574
$ws.Range("H29").Value = 6437.5
$ws.Range("I29").Value = 1750
$ws.Range("J29").Value = 8000
$ws.Range("K29").Value = 1750
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = -1457
$ws.Range("N29").Value = -8586
$ws.Range("H31").Value = 2590.175
$ws.Range("I31").Value = 1866.375
$ws.Range("J31").Value = 5485.375
$ws.Range("K31").Value = 1866.375
$ws.Range("L31").Value = 5485.375
$ws.Range("M31").Value = -1571.375
$ws.Range("N31").Value = -6075.375
$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -2632
$ws.Range("H34").Value = 2590.175
$ws.Range("I34").Value = 1866.375
$ws.Range("J34").Value = 5485.375
$ws.Range("K34").Value = 1866.375
$ws.Range("L34").Value = 5485.375
$ws.Range("M34").Value = -1664.375
$ws.Range("N34").Value = -5889.375
$ws.Range("H35").Value = 1107.1428
$ws.Range("I35").Value = 791.6667
$ws.Range("K35").Value = 791.6667
$ws.Range("M35").Value = -497.6667
$ws.Range("H36").Value = 1750
$ws.Range("I36").Value = 1000
$ws.Range("J36").Value = 2500
$ws.Range("K36").Value = 1000
$ws.Range("L36").Value = 2500
$ws.Range("M36").Value = -612
$ws.Range("N36").Value = -3276
$ws.Range("H40").Value = 1750
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -840
$ws.Range("N40").Value = -2820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 59184190
$ws.Range("H122").Value = 1084.3043
$ws.Range("I122").Value = 891.8333
$ws.Range("J122").Value = 1152.2354
$ws.Range("K122").Value = 8026.4997
$ws.Range("L122").Value = 10370.1186
$ws.Range("M122").Value = -5576.4997
$ws.Range("N122").Value = -15270.1186

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 51493.523
$ws.Range("I100").Value = 223665.8
$ws.Range("J100").Value = 3667.889
$ws.Range("K100").Value = 223665.8
$ws.Range("L100").Value = 3667.889
$ws.Range("M100").Value = -223124.8
$ws.Range("N100").Value = -4749.889
$ws.Range("H122").Value = 5672.64
$ws.Range("I122").Value = 2560.5
$ws.Range("K122").Value = 7681.5
$ws.Range("M122").Value = -5231.5
$ws.Range("H136").Value = 2803.17
$ws.Range("I136").Value = 1540.5758
$ws.Range("J136").Value = 4886.45
$ws.Range("K136").Value = 4621.7274
$ws.Range("L136").Value = 14659.35
$ws.Range("M136").Value = -2071.7274
$ws.Range("N136").Value = -19759.35

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 236796.08
$ws.Range("I122").Value = 2115.8684
$ws.Range("J122").Value = 1128580.9
$ws.Range("K122").Value = 6347.6052
$ws.Range("L122").Value = 3385742.7
$ws.Range("M122").Value = -3897.6052
$ws.Range("N122").Value = -3390642.7
$ws.Range("H132").Value = 1848.1351
$ws.Range("I132").Value = 1297.4828
$ws.Range("J132").Value = 3844.25
$ws.Range("K132").Value = 3892.4484
$ws.Range("L132").Value = 11532.75
$ws.Range("M132").Value = -1362.4484
$ws.Range("N132").Value = -16592.75
